$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Update row 8 data: date moves from 09-Mar-22 to 16-Mar-22, and the
# mobility indicator percentages are refreshed for the new week.
$ws.Range("A8").Value = 44636
$ws.Range("B8").Value = 0.1
$ws.Range("C8").Value = 0.26
$ws.Range("D8").Value = 0.21
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = 0.12

# Update the view to match where the user left the selection/scroll.
$ws.Range("E8").Select()
$excel.ActiveWindow.ScrollRow = 11
